$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-12 with new data
$ws.Range("A2").Value = "Department of Energy"
$ws.Range("B2").Value = "Caroline Grey"
$ws.Range("E2").Value = "Biden for President"

$ws.Range("A3").Value = "Department of Health and Human Services"
$ws.Range("B3").Value = "Clara Pratte"
$ws.Range("D3").Value = "Y"
$ws.Range("E3").Value = "Strongbow Strategies"

$ws.Range("A4").Value = "Department of Housing and Urban Development"
$ws.Range("B4").Value = "Analysse Escobar"
$ws.Range("E4").Value = "Biden for President"
$ws.Range("F4").Value = "Transition — PT Fund, Inc."

$ws.Range("A5").Value = "Department of Housing and Urban Development"
$ws.Range("B5").Value = "Mikayla Ferrell"
$ws.Range("E5").Value = "Pennsylvania Democratic Party"
$ws.Range("F5").Value = "Transition — PT Fund, Inc."

$ws.Range("A6").Value = "Department of Justice"
$ws.Range("B6").Value = "Theresa Bradley"
$ws.Range("E6").Value = "Biden for President"
$ws.Range("F6").Value = "Transition — PT Fund, Inc."

$ws.Range("A7").Value = "Department of the Interior"
$ws.Range("B7").Value = "Maggie Thompson"
$ws.Range("E7").Value = "Biden for President"
$ws.Range("F7").Value = "Transition — PT Fund, Inc."

$ws.Range("A8").Value = "Department of the Treasury"
$ws.Range("B8").Value = "William Doerrer"
$ws.Range("E8").Value = "Biden for President"

$ws.Range("A9").Value = "Department of Transportation"
$ws.Range("B9").Value = "Allie Panther"
$ws.Range("E9").Value = "Biden for President"
$ws.Range("F9").Value = "Transition — PT Fund, Inc."

$ws.Range("A10").Value = "Environmental Protection Agency"
$ws.Range("B10").Value = "Sinceré Harris"
$ws.Range("E10").Value = "Pennsylvania Democratic Party"
$ws.Range("F10").Value = "Transition — PT Fund, Inc."

$ws.Range("A11").Value = "Office of Personnel Management"
$ws.Range("B11").Value = "Jason Tengco"
$ws.Range("E11").Value = "Biden for President"
$ws.Range("F11").Value = "Transition — PT Fund, Inc."

$ws.Range("A12").Value = "United States Department of Agriculture"
$ws.Range("B12").Value = "Lexi Coburn"
$ws.Range("E12").Value = "North Carolina Democratic Party"

# Delete old row 13 (its data has been folded into row 12 shift; remove entirely)
$ws.Rows.Item(13).Delete()
